$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-19 18:34:59"
$wsZhCn.Range("H4").Value = "2016-03-19 18:35:20"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-19 18:35:03"
$wsDeDe.Range("H4").Value = "2016-03-19 18:35:25"
